# 1401EE17.xlsx fixes: zero out the student's right/wrong counts for the
# first section, mark the second section as "Absent", and clear out the
# stray "Student Ans" values that shouldn't have been populated
# (fixes related to fileNames sent in emails / re-grading).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Score summary block (rows 10 & 12) ---------------------------------
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 28

$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = "Absent"

# --- Clear stray "Student Ans" entries (columns A & D) ------------------
# These rows incorrectly had an answer recorded; they should be blank,
# matching the style used by the other blank "Student Ans" cells.
$rowsColA = @(16, 18, 21, 22, 26, 28, 30, 31, 32, 33, 35, 36, 37, 38, 39)
foreach ($r in $rowsColA) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.ClearContents()
    $cell.Style = "normalStyle"
}

$rowsColD = @(16, 18)
foreach ($r in $rowsColD) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.ClearContents()
    $cell.Style = "normalStyle"
}
